$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 6, pushing every
# existing record (previously rows 6-130) down by one row.
$ws.Rows("6").Insert()

# Seed the new row with the same "constant" columns as its neighbour (now
# row 7, the record that used to live in row 6) so formatting/style (e.g.
# the date number format on column D) and the repeated descriptive fields
# come along for free, then overwrite the fields that actually carry the
# new weekly observation.
$ws.Range("A7:R7").Copy()
$ws.Range("A6:R6").PasteSpecial()
$excel.CutCopyMode = $false

$ws.Range("D6").Value2 = 45190
$ws.Range("J6").Value2 = 350
$ws.Range("K6").Value2 = 2200
$ws.Range("L6").Value2 = 2500
$ws.Range("M6").Value2 = 2371
$ws.Range("P6").Value2 = 1186
